# cryptos.xlsx -- refresh "Price" (D) and "Volume(1h)" (E) columns with
# the latest scraped coin snapshot (GitHub Actions cron update).
#
# The Price column holds values that read like numbers (e.g. "1.638.41",
# "15.60", "0.0₅ 7748") but are stored as literal text in the workbook --
# some aren't valid numbers at all (two decimal points / subscript digits),
# and the ones that are would silently lose significant trailing zeros if
# Excel auto-coerced them to numeric (e.g. "15.60" -> 15.6). Force each
# Price cell to Text before writing, then drop the number-format override
# so the change doesn't leave a stray style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '25.905.55'
$dCell.ClearFormats()
$ws.Range("E2").Value = '  -0.10%  '
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.638.41'
$dCell.ClearFormats()
$ws.Range("E3").Value = '  -0.15%  '
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '1.002'
$dCell.ClearFormats()
$ws.Range("E4").Value = '  -0.28%  '
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '214.63'
$dCell.ClearFormats()
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("E7").Value = '  -0.24%  '
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '0.2575'
$dCell.ClearFormats()
$ws.Range("E8").Value = '  +0.32%  '
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.06363'
$dCell.ClearFormats()
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("E10").Value = '  +1.74%  '
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.07742'
$dCell.ClearFormats()
$ws.Range("E11").Value = '  -0.36%  '
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '4.297'
$dCell.ClearFormats()
$ws.Range("E12").Value = '  -0.03%  '
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '1.640.65'
$dCell.ClearFormats()
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("E14").Value = '  +0.27%  '
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₅7748'
$dCell.ClearFormats()
$ws.Range("E15").Value = '  -1.56%  '
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '64.17'
$dCell.ClearFormats()
$ws.Range("E16").Value = '  -0.55%  '
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '25.941.32'
$dCell.ClearFormats()
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("E18").Value = '  -0.30%  '
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '4.465'
$dCell.ClearFormats()
$ws.Range("E19").Value = '  +0.80%  '
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '195.79'
$dCell.ClearFormats()
$ws.Range("E20").Value = '  -0.89%  '
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '9.952'
$dCell.ClearFormats()
$ws.Range("E21").Value = '  -0.03%  '
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '6.145'
$dCell.ClearFormats()
$ws.Range("E22").Value = '  +1.68%  '
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '1.002'
$dCell.ClearFormats()
$ws.Range("E23").Value = '  -0.39%  '
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '1.899'
$dCell.ClearFormats()
$ws.Range("E24").Value = '  +0.94%  '
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '143.35'
$dCell.ClearFormats()
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("E26").Value = '  +10.22%  '
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '6.849'
$dCell.ClearFormats()
$ws.Range("E27").Value = '  -0.43%  '
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '15.60'
$dCell.ClearFormats()
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -2.81%  '
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '3.251'
$dCell.ClearFormats()
$ws.Range("E31").Value = '  -0.45%  '
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '3.206'
$dCell.ClearFormats()
$ws.Range("E32").Value = '  +0.43%  '
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '1.557'
$dCell.ClearFormats()
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("E34").Value = '  +0.50%  '
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9160'
$dCell.ClearFormats()
$ws.Range("E35").Value = '  +2.48%  '
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '2.569'
$dCell.ClearFormats()
$ws.Range("E36").Value = '  -1.15%  '
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '1.137.55'
$dCell.ClearFormats()
$ws.Range("E37").Value = '  +0.25%  '
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5524'
$dCell.ClearFormats()
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("E39").Value = '  +0.74%  '
$ws.Range("E40").Value = '  -0.28%  '
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '5.598'
$dCell.ClearFormats()
$ws.Range("E41").Value = '  -0.61%  '
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '0.8043'
$dCell.ClearFormats()
$ws.Range("E42").Value = '  -1.30%  '
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '98.46'
$dCell.ClearFormats()
$ws.Range("E43").Value = '  -1.49%  '
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₈122'
$dCell.ClearFormats()
$ws.Range("E44").Value = '  -8.77%  '
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '1.771.34'
$dCell.ClearFormats()
$ws.Range("E45").Value = '  -0.48%  '
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4494'
$dCell.ClearFormats()
$ws.Range("E46").Value = '  -0.76%  '
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '55.24'
$dCell.ClearFormats()
$ws.Range("E47").Value = '  +0.73%  '
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.ClearFormats()
$ws.Range("E48").Value = '  -0.37%  '
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '0.05181'
$dCell.ClearFormats()
$ws.Range("E49").Value = '  +2.10%  '
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '7.534'
$dCell.ClearFormats()
$ws.Range("E50").Value = '  +1.66%  '
$ws.Range("E51").Value = '  -0.34%  '
